$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Summon Rating (column D) values for rows with recalculated odds ---
$ws.Range("D2:D94").NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "14.133340412625804"
$ws.Cells.Item(3, 4).Value = "19.077895625846853"
$ws.Cells.Item(4, 4).Value = "16.868613330468833"
$ws.Cells.Item(5, 4).Value = "10.352530283634735"
$ws.Cells.Item(6, 4).Value = "10.34971053072163"
$ws.Cells.Item(8, 4).Value = "6.510814244639059"
$ws.Cells.Item(9, 4).Value = "5.578333906768179"
$ws.Cells.Item(13, 4).Value = "7.309310127153189"
$ws.Cells.Item(14, 4).Value = "10.940510380601005"
$ws.Cells.Item(15, 4).Value = "2.6268567181014317"
$ws.Cells.Item(16, 4).Value = "22.2859941312868"
$ws.Cells.Item(18, 4).Value = "9.389287054125319"
$ws.Cells.Item(19, 4).Value = "14.190860551586304"
$ws.Cells.Item(20, 4).Value = "10.98307009639102"
$ws.Cells.Item(23, 4).Value = "4.082528974969949"
$ws.Cells.Item(24, 4).Value = "27.08002772244592"
$ws.Cells.Item(25, 4).Value = "21.97372854706214"
$ws.Cells.Item(26, 4).Value = "4.462902158067706"
$ws.Cells.Item(27, 4).Value = "7.114561657913072"
$ws.Cells.Item(29, 4).Value = "15.724323622620199"
$ws.Cells.Item(30, 4).Value = "8.877955651729955"
$ws.Cells.Item(31, 4).Value = "2.685810565637707"
$ws.Cells.Item(32, 4).Value = "3.067226924536795"
$ws.Cells.Item(33, 4).Value = "6.59582218088197"
$ws.Cells.Item(36, 4).Value = "64.78687639641227"
$ws.Cells.Item(37, 4).Value = "14.134038143912903"
$ws.Cells.Item(38, 4).Value = "13.745754425097564"
$ws.Cells.Item(40, 4).Value = "2.83537508308558"
$ws.Cells.Item(42, 4).Value = "2.5284272519392257"
$ws.Cells.Item(43, 4).Value = "62.34850144228113"
$ws.Cells.Item(44, 4).Value = "6.795075598727976"
$ws.Cells.Item(45, 4).Value = "11.34255459010593"
$ws.Cells.Item(46, 4).Value = "1.0477672544943877"
$ws.Cells.Item(47, 4).Value = "3.2532313569149847"
$ws.Cells.Item(49, 4).Value = "3.2795068569092893"
$ws.Cells.Item(50, 4).Value = "2.1615896714057374"
$ws.Cells.Item(51, 4).Value = "6.865860983112577"
$ws.Cells.Item(52, 4).Value = "2.921017307124261"
$ws.Cells.Item(53, 4).Value = "1.2962962962962958"
$ws.Cells.Item(54, 4).Value = "20.311560476932023"
$ws.Cells.Item(56, 4).Value = "4.669913530351197"
$ws.Cells.Item(58, 4).Value = "5.655445028056643"
$ws.Cells.Item(59, 4).Value = "4.900064301327429"
$ws.Cells.Item(62, 4).Value = "4.117636032410978"
$ws.Cells.Item(63, 4).Value = "3.4799681180262105"
$ws.Cells.Item(64, 4).Value = "12.612343854748275"
$ws.Cells.Item(65, 4).Value = "12.18306577064891"
$ws.Cells.Item(66, 4).Value = "16.07008016401693"
$ws.Cells.Item(67, 4).Value = "83.1312508532607"
$ws.Cells.Item(68, 4).Value = "96.55930484816379"
$ws.Cells.Item(69, 4).Value = "2.3591775203152823"
$ws.Cells.Item(70, 4).Value = "3.831163045718534"
$ws.Cells.Item(71, 4).Value = "2.383002484969857"
$ws.Cells.Item(72, 4).Value = "3.6776361212976343"
$ws.Cells.Item(73, 4).Value = "12.36190870176738"
$ws.Cells.Item(74, 4).Value = "4.128562334741663"
$ws.Cells.Item(75, 4).Value = "1.386542352415373"
$ws.Cells.Item(76, 4).Value = "4.836681769311175"
$ws.Cells.Item(77, 4).Value = "11.908048677191992"
$ws.Cells.Item(78, 4).Value = "2.974753956948206"
$ws.Cells.Item(79, 4).Value = "1.9379151650275726"
$ws.Cells.Item(80, 4).Value = "39.31061310425023"
$ws.Cells.Item(81, 4).Value = "8.179683151192107"
$ws.Cells.Item(82, 4).Value = "98.09400349519292"
$ws.Cells.Item(83, 4).Value = "2.271869476881736"
$ws.Cells.Item(84, 4).Value = "52.51537479829067"
$ws.Cells.Item(85, 4).Value = "24.76966023485272"
$ws.Cells.Item(86, 4).Value = "13.803210325440759"
$ws.Cells.Item(87, 4).Value = "4.659542244013985"
$ws.Cells.Item(88, 4).Value = "1.4287527744381399"
$ws.Cells.Item(89, 4).Value = "2.7266508617553966"
$ws.Cells.Item(90, 4).Value = "0.7004388124634752"
$ws.Cells.Item(94, 4).Value = "11.798349300571095"
$ws.Range("D2:D94").Style = "Normal"

# --- Append newly summoned/unlocked units as additional rows ---
$ws.Range("A95:D103").NumberFormat = "@"
$ws.Cells.Item(95, 1).Value = "94"
$ws.Cells.Item(95, 2).Value = "DFLR_PHY_Buutenks_"
$ws.Cells.Item(95, 3).Value = "5"
$ws.Cells.Item(95, 4).Value = "0.2"
$ws.Cells.Item(96, 1).Value = "95"
$ws.Cells.Item(96, 2).Value = "DFLR_STR_SS_Vegito"
$ws.Cells.Item(96, 3).Value = "5"
$ws.Cells.Item(96, 4).Value = "0.2"
$ws.Cells.Item(97, 1).Value = "96"
$ws.Cells.Item(97, 2).Value = "BU_INT_Dodoria_"
$ws.Cells.Item(97, 3).Value = "0"
$ws.Cells.Item(97, 4).Value = "16.43953697967035"
$ws.Cells.Item(98, 1).Value = "97"
$ws.Cells.Item(98, 2).Value = "DF_AGL_1stForm_Frieza"
$ws.Cells.Item(98, 3).Value = "0"
$ws.Cells.Item(98, 4).Value = "48.09467853489954"
$ws.Cells.Item(99, 1).Value = "98"
$ws.Cells.Item(99, 2).Value = "DF_AGL_SS_Bardock"
$ws.Cells.Item(99, 3).Value = "5"
$ws.Cells.Item(99, 4).Value = "0.2"
$ws.Cells.Item(100, 1).Value = "99"
$ws.Cells.Item(100, 2).Value = "BU_AGL_King_Vegeta"
$ws.Cells.Item(100, 3).Value = "0"
$ws.Cells.Item(100, 4).Value = "12.20756811006764"
$ws.Cells.Item(101, 1).Value = "100"
$ws.Cells.Item(101, 2).Value = "DF_TEQ_SS_Vegeta"
$ws.Cells.Item(101, 3).Value = "0"
$ws.Cells.Item(101, 4).Value = "75.73779918278534"
$ws.Cells.Item(102, 1).Value = "101"
$ws.Cells.Item(102, 2).Value = "DF_AGL_DragonFist_Goku"
$ws.Cells.Item(102, 3).Value = "2"
$ws.Cells.Item(102, 4).Value = "4.0054013442748255"
$ws.Cells.Item(103, 1).Value = "102"
$ws.Cells.Item(103, 2).Value = "DFLR_INT_SFPS4_Goku"
$ws.Cells.Item(103, 3).Value = "1"
$ws.Cells.Item(103, 4).Value = "5.968161313573148"
$ws.Range("A95:D103").Style = "Normal"

# --- Re-apply the bold/centered/bordered ID-column format to the new rows (copy from an existing ID cell) ---
$ws.Range("A2").Copy()
$ws.Range("A95:A103").PasteSpecial(-4122)
$excel.CutCopyMode = 0
